$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Round Title" column (AF) is being removed entirely; everything to its
# right shifts one column to the left (AG->AF, AH->AG, ... AW->AV).
$ws.Columns("AF:AF").Delete()

# Keep the hidden AutoFilter-database defined name in sync with the new,
# one-column-narrower used range.
$n = $wb.Names.Item("Function!_FilterDatabase")
$n.RefersTo = '=Function!$A$1:$AV$2'

# "Updated Code Left split": the Code column header (A1) now carries an
# explicit Text number format (in addition to its existing left alignment).
$ws.Range("A1").NumberFormat = "@"

Write-Output "done"
